$d = $word.ActiveDocument

# 1. Merge "A05:2021-Security Misconfigura" + "t" + "ion" into a single run's text.
#    (Text content is unchanged overall, just de-fragmenting runs -> simple text replace is fine.)
$d.Content.Find.Execute("A05:2021-Security Misconfiguration", $false, $false, $false, $false, $false, $true, 1, $false, "A05:2021-Security Misconfiguration", 2) | Out-Null

# 2. Add " Could use OAuth." after "...by users with correct permissions."
$d.Content.Find.Execute("by users with correct permissions.", $true, $false, $false, $false, $false, $true, 1, $false, "by users with correct permissions. Could use OAuth.", 2) | Out-Null

# 3. Replace "Use of Auth0" with "N/A"
$d.Content.Find.Execute("Use of Auth0", $true, $false, $false, $false, $false, $true, 1, $false, "N/A", 2) | Out-Null

# 4. Remove trailing period and add " Through email login links for auth."
$d.Content.Find.Execute("could be used to counteract this vulnerability.", $true, $false, $false, $false, $false, $true, 1, $false, "could be used to counteract this vulnerability Through email login links for auth.", 2) | Out-Null

# 5. Merge "A09:2021-Security Logg" + "i" + "ng and Monitoring Failures" into a single run's text.
$d.Content.Find.Execute("A09:2021-Security Logging and Monitoring Failures", $false, $false, $false, $false, $false, $true, 1, $false, "A09:2021-Security Logging and Monitoring Failures", 2) | Out-Null
